$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.146235
$ws.Range("H2").Value = 0.438705
$ws.Range("I2").Value = 0.0224838618501081
$ws.Range("J2").Value = 0.0224838618501081
$ws.Range("M2").Value = 1.594873333333333
$ws.Range("N2").Value = 4.78462
$ws.Range("O2").Value = 0.09372679355272211
$ws.Range("P2").Value = 0.09372679355272213
$ws.Range("Q2").Value = 0.2332263019
$ws.Range("R2").Value = 2.0990367171
$ws.Range("S2").Value = 0.002107340277893006
$ws.Range("T2").Value = 0.002107340277893007
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.146235
$ws.Range("H3").Value = 0.438705
$ws.Range("I3").Value = 0.0224838618501081
$ws.Range("J3").Value = 0.0224838618501081
$ws.Range("O3").Value = 0.2690834924840127
$ws.Range("P3").Value = 0.2690834924840128
$ws.Range("Q3").Value = 0.669577454595
$ws.Range("R3").Value = 6.026197091355
$ws.Range("S3").Value = 0.006050036071155142
$ws.Range("T3").Value = 0.006050036071155145
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.146235
$ws.Range("H4").Value = 0.438705
$ws.Range("I4").Value = 0.0224838618501081
$ws.Range("J4").Value = 0.0224838618501081
$ws.Range("M4").Value = 4.495828
$ws.Range("N4").Value = 13.487484
$ws.Range("O4").Value = 0.2642087832291055
$ws.Range("P4").Value = 0.2642087832291055
$ws.Range("Q4").Value = 0.6574474075800001
$ws.Range("R4").Value = 5.91702666822
$ws.Range("S4").Value = 0.005940433781708366
$ws.Range("T4").Value = 0.005940433781708367
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.146235
$ws.Range("H5").Value = 0.438705
$ws.Range("I5").Value = 0.0224838618501081
$ws.Range("J5").Value = 0.0224838618501081
$ws.Range("M5").Value = 6.346716
$ws.Range("N5").Value = 19.040148
$ws.Range("O5").Value = 0.3729809307341596
$ws.Range("P5").Value = 0.3729809307341597
$ws.Range("Q5").Value = 0.92811201426
$ws.Range("R5").Value = 8.353008128339999
$ws.Range("S5").Value = 0.008386051719351582
$ws.Range("T5").Value = 0.008386051719351586
$ws.Range("G6").Value = 3.793107666666666
$ws.Range("I6").Value = 0.5831962851568996
$ws.Range("J6").Value = 0.5831962851568997
$ws.Range("M6").Value = 1.594873333333333
$ws.Range("N6").Value = 4.78462
$ws.Range("O6").Value = 0.09372679355272211
$ws.Range("P6").Value = 0.09372679355272213
$ws.Range("Q6").Value = 6.049526268028888
$ws.Range("R6").Value = 54.44573641226
$ws.Range("S6").Value = 0.05466111781961518
$ws.Range("T6").Value = 0.0546611178196152
$ws.Range("G7").Value = 3.793107666666666
$ws.Range("I7").Value = 0.5831962851568996
$ws.Range("J7").Value = 0.5831962851568997
$ws.Range("O7").Value = 0.2690834924840127
$ws.Range("P7").Value = 0.2690834924840128
$ws.Range("S7").Value = 0.1569284932137207
$ws.Range("T7").Value = 0.1569284932137208
$ws.Range("G8").Value = 3.793107666666666
$ws.Range("I8").Value = 0.5831962851568996
$ws.Range("J8").Value = 0.5831962851568997
$ws.Range("M8").Value = 4.495828
$ws.Range("N8").Value = 13.487484
$ws.Range("O8").Value = 0.2642087832291055
$ws.Range("P8").Value = 0.2642087832291055
$ws.Range("Q8").Value = 17.05315965481467
$ws.Range("R8").Value = 153.478436893332
$ws.Range("S8").Value = 0.1540855808850389
$ws.Range("T8").Value = 0.1540855808850389
$ws.Range("G9").Value = 3.793107666666666
$ws.Range("I9").Value = 0.5831962851568996
$ws.Range("J9").Value = 0.5831962851568997
$ws.Range("M9").Value = 6.346716
$ws.Range("N9").Value = 19.040148
$ws.Range("O9").Value = 0.3729809307341596
$ws.Range("P9").Value = 0.3729809307341597
$ws.Range("Q9").Value = 24.073777117756
$ws.Range("R9").Value = 216.663994059804
$ws.Range("S9").Value = 0.2175210932385248
$ws.Range("T9").Value = 0.2175210932385248
$ws.Range("G10").Value = 2.288493
$ws.Range("H10").Value = 6.865479000000001
$ws.Range("I10").Value = 0.3518594075080483
$ws.Range("J10").Value = 0.3518594075080483
$ws.Range("M10").Value = 1.594873333333333
$ws.Range("N10").Value = 4.78462
$ws.Range("O10").Value = 0.09372679355272211
$ws.Range("P10").Value = 0.09372679355272213
$ws.Range("Q10").Value = 3.64985645922
$ws.Range("R10").Value = 32.84870813298001
$ws.Range("S10").Value = 0.03297865404708996
$ws.Range("T10").Value = 0.03297865404708997
$ws.Range("G11").Value = 2.288493
$ws.Range("H11").Value = 6.865479000000001
$ws.Range("I11").Value = 0.3518594075080483
$ws.Range("J11").Value = 0.3518594075080483
$ws.Range("O11").Value = 0.2690834924840127
$ws.Range("P11").Value = 0.2690834924840128
$ws.Range("Q11").Value = 10.478499113061
$ws.Range("R11").Value = 94.306492017549
$ws.Range("S11").Value = 0.09467955823562108
$ws.Range("T11").Value = 0.09467955823562109
$ws.Range("G12").Value = 2.288493
$ws.Range("H12").Value = 6.865479000000001
$ws.Range("I12").Value = 0.3518594075080483
$ws.Range("J12").Value = 0.3518594075080483
$ws.Range("M12").Value = 4.495828
$ws.Range("N12").Value = 13.487484
$ws.Range("O12").Value = 0.2642087832291055
$ws.Range("P12").Value = 0.2642087832291055
$ws.Range("Q12").Value = 10.288670907204
$ws.Range("R12").Value = 92.59803816483601
$ws.Range("S12").Value = 0.09296434592541543
$ws.Range("T12").Value = 0.09296434592541543
$ws.Range("G13").Value = 2.288493
$ws.Range("H13").Value = 6.865479000000001
$ws.Range("I13").Value = 0.3518594075080483
$ws.Range("J13").Value = 0.3518594075080483
$ws.Range("M13").Value = 6.346716
$ws.Range("N13").Value = 19.040148
$ws.Range("O13").Value = 0.3729809307341596
$ws.Range("P13").Value = 0.3729809307341597
$ws.Range("Q13").Value = 14.524415138988
$ws.Range("R13").Value = 130.719736250892
$ws.Range("S13").Value = 0.1312368492999218
$ws.Range("T13").Value = 0.1312368492999218
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.2761626666666666
$ws.Range("H14").Value = 0.8284879999999999
$ws.Range("I14").Value = 0.04246044548494399
$ws.Range("J14").Value = 0.042460445484944
$ws.Range("M14").Value = 1.594873333333333
$ws.Range("N14").Value = 4.78462
$ws.Range("O14").Value = 0.09372679355272211
$ws.Range("P14").Value = 0.09372679355272213
$ws.Range("Q14").Value = 0.4404444727288888
$ws.Range("R14").Value = 3.96400025456
$ws.Range("S14").Value = 0.003979681408123957
$ws.Range("T14").Value = 0.003979681408123959
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.2761626666666666
$ws.Range("H15").Value = 0.8284879999999999
$ws.Range("I15").Value = 0.04246044548494399
$ws.Range("J15").Value = 0.042460445484944
$ws.Range("O15").Value = 0.2690834924840127
$ws.Range("P15").Value = 0.2690834924840128
$ws.Range("Q15").Value = 1.264487266392
$ws.Range("R15").Value = 11.380385397528
$ws.Range("S15").Value = 0.01142540496351576
$ws.Range("T15").Value = 0.01142540496351576
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.2761626666666666
$ws.Range("H16").Value = 0.8284879999999999
$ws.Range("I16").Value = 0.04246044548494399
$ws.Range("J16").Value = 0.042460445484944
$ws.Range("M16").Value = 4.495828
$ws.Range("N16").Value = 13.487484
$ws.Range("O16").Value = 0.2642087832291055
$ws.Range("P16").Value = 0.2642087832291055
$ws.Range("Q16").Value = 1.241579849354667
$ws.Range("R16").Value = 11.174218644192
$ws.Range("S16").Value = 0.01121842263694282
$ws.Range("T16").Value = 0.01121842263694282
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.2761626666666666
$ws.Range("H17").Value = 0.8284879999999999
$ws.Range("I17").Value = 0.04246044548494399
$ws.Range("J17").Value = 0.042460445484944
$ws.Range("M17").Value = 6.346716
$ws.Range("N17").Value = 19.040148
$ws.Range("O17").Value = 0.3729809307341596
$ws.Range("P17").Value = 0.3729809307341597
$ws.Range("Q17").Value = 1.752726015136
$ws.Range("R17").Value = 15.774534136224
$ws.Range("S17").Value = 0.01583693647636146
$ws.Range("T17").Value = 0.01583693647636146
